# Apply the KSA_Cities.xlsx update:
#  1. Correct two mis-spelled Arabic city names already present in the sheet.
#  2. Append three new city rows at the bottom of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Fix existing Arabic spellings -------------------------------------
# Row 44  : Dheba / Duba, Saudi Arabia  -> Arabic name "ضبا" becomes "ضباء"
$ws.Cells.Item(44, 3).Value = "ضباء"

# Row 117 : Al Suhan                    -> Arabic name "الصحن" becomes "السحن"
$ws.Cells.Item(117, 3).Value = "السحن"

# --- 2. Append new rows -----------------------------------------------------
# Row 161: Ranyah
$ws.Cells.Item(161, 1).Value = "Ranyah"
$ws.Cells.Item(161, 2).Value = "Ranyah"
$ws.Cells.Item(161, 3).Value = "رنية"
$ws.Cells.Item(161, 4).Value = 21.263856000000001
$ws.Cells.Item(161, 5).Value = 42.853374000000002
$ws.Cells.Item(161, 6).Value = "منطقة مكة المكرمة"
$ws.Cells.Item(161, 7).Value = "غرب المملكة"

# Row 162: Dhamad
$ws.Cells.Item(162, 1).Value = "Dhamad"
$ws.Cells.Item(162, 2).Value = "Dhamad"
$ws.Cells.Item(162, 3).Value = "ضمد"
$ws.Cells.Item(162, 4).Value = 17.110188999999998
$ws.Cells.Item(162, 5).Value = 42.775911000000001
$ws.Cells.Item(162, 6).Value = "منطقة جازان"
$ws.Cells.Item(162, 7).Value = "جنوب المملكة"

# Row 163: Nafy
$ws.Cells.Item(163, 1).Value = "Nafy"
$ws.Cells.Item(163, 2).Value = "Nafy"
$ws.Cells.Item(163, 3).Value = "نفى"
$ws.Cells.Item(163, 4).Value = 25.014095000000001
$ws.Cells.Item(163, 5).Value = 43.817810000000001
$ws.Cells.Item(163, 6).Value = "منطقة الرياض"
$ws.Cells.Item(163, 7).Value = "وسط المملكة"
